$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Treatment"), shifting Treatment/Analytes/append-note
# one column to the right so the new "Subject id" column can be added.
$ws.Columns("D:D").Insert()

# New header cell: "Subject id" in bold black Verdana, centered (matches the look of the
# other header font but as its own distinct style).
$ws.Range("D1").Value = "Subject id"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Name = "Verdana"
$ws.Range("D1").Font.Size = 11
$ws.Range("D1").Font.Color = 0
$ws.Range("D1").HorizontalAlignment = -4108

# The sample id value (7101) used to live in B2; it now belongs in the new D2 cell.
# Copy B2's current (pre-edit) number format/border/alignment into D2 first, then it
# already carries the right value too.
$ws.Range("B2").Copy($ws.Range("D2"))

# B2 becomes the new "Subject id" value cell - plain left-aligned Verdana text, no border.
$ws.Range("B2").Value2 = "MU_001"
$ws.Range("B2").Font.Name = "Verdana"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").Borders.LineStyle = 0

# Selection, as recorded in the edited workbook.
$ws.Range("B6").Select()
